$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1843971631205674
$ws.Range("C2").Value = 0.5815602836879432
$ws.Range("J2").Value = 0.01418439716312057
$ws.Range("P2").Value = 0.1099290780141844
$ws.Range("S2").Value = 0.1099290780141844
$ws.Range("B3").Value = 0.006060606060606061
$ws.Range("C3").Value = 0.02424242424242424
$ws.Range("J3").Value = 0.0303030303030303
$ws.Range("P3").Value = 0.7333333333333333
$ws.Range("S3").Value = 0.2060606060606061
$ws.Range("J4").Value = 0.1111111111111111
$ws.Range("P4").Value = 0.6222222222222222
$ws.Range("S4").Value = 0.2666666666666667
$ws.Range("B6").Value = 0.05836575875486381
$ws.Range("D6").Value = 0.01556420233463035
$ws.Range("F6").Value = 0.05058365758754864
$ws.Range("J6").Value = 0.3035019455252918
$ws.Range("O6").Value = 0.01556420233463035
$ws.Range("Q6").Value = 0.1556420233463035
$ws.Range("R6").Value = 0.07782101167315175
$ws.Range("S6").Value = 0.3229571984435798
$ws.Range("B7").Value = 0.06593406593406594
$ws.Range("D7").Value = 0.02747252747252747
$ws.Range("F7").Value = 0.04945054945054945
$ws.Range("J7").Value = 0.1208791208791209
$ws.Range("O7").Value = 0.02197802197802198
$ws.Range("Q7").Value = 0.2197802197802198
$ws.Range("R7").Value = 0.09340659340659341
$ws.Range("S7").Value = 0.4010989010989011
$ws.Range("B8").Value = 0.09808612440191387
$ws.Range("D8").Value = 0.0215311004784689
$ws.Range("F8").Value = 0.07655502392344497
$ws.Range("J8").Value = 0.0861244019138756
$ws.Range("O8").Value = 0.01674641148325359
$ws.Range("Q8").Value = 0.1842105263157895
$ws.Range("R8").Value = 0.1124401913875598
$ws.Range("S8").Value = 0.4043062200956938
$ws.Range("B9").Value = 0.08415841584158416
$ws.Range("D9").Value = 0.009900990099009901
$ws.Range("F9").Value = 0.05445544554455446
$ws.Range("J9").Value = 0.07425742574257425
$ws.Range("O9").Value = 0.0297029702970297
$ws.Range("Q9").Value = 0.1930693069306931
$ws.Range("R9").Value = 0.07920792079207921
$ws.Range("S9").Value = 0.4752475247524752
$ws.Range("B10").Value = 0.1180781758957655
$ws.Range("D10").Value = 0.02198697068403909
$ws.Range("F10").Value = 0.07491856677524431
$ws.Range("J10").Value = 0.0985342019543974
$ws.Range("O10").Value = 0.01710097719869707
$ws.Range("Q10").Value = 0.2182410423452769
$ws.Range("R10").Value = 0.08794788273615635
$ws.Range("S10").Value = 0.3631921824104234
$ws.Range("G11").Value = 0.1595744680851064
$ws.Range("J11").Value = 0.0851063829787234
$ws.Range("K11").Value = 0.2127659574468085
$ws.Range("L11").Value = 0.5141843971631206
$ws.Range("S11").Value = 0.02836879432624113
$ws.Range("G12").Value = 0.7337662337662337
$ws.Range("J12").Value = 0.2142857142857143
$ws.Range("L12").Value = 0.03896103896103896
$ws.Range("S12").Value = 0.01298701298701299
$ws.Range("F13").Value = 0.02325581395348837
$ws.Range("G13").Value = 0.6744186046511628
$ws.Range("J13").Value = 0.2790697674418605
$ws.Range("S13").Value = 0.02325581395348837
$ws.Range("F15").Value = 0.03557312252964427
$ws.Range("H15").Value = 0.150197628458498
$ws.Range("I15").Value = 0.06324110671936758
$ws.Range("J15").Value = 0.3715415019762846
$ws.Range("K15").Value = 0.05138339920948617
$ws.Range("M15").Value = 0.01185770750988142
$ws.Range("O15").Value = 0.06719367588932806
$ws.Range("S15").Value = 0.2490118577075099
$ws.Range("F16").Value = 0.02312138728323699
$ws.Range("H16").Value = 0.1098265895953757
$ws.Range("I16").Value = 0.09248554913294797
$ws.Range("J16").Value = 0.3988439306358382
$ws.Range("K16").Value = 0.1040462427745665
$ws.Range("M16").Value = 0.02312138728323699
$ws.Range("O16").Value = 0.09826589595375723
$ws.Range("S16").Value = 0.1502890173410405
$ws.Range("F17").Value = 0.03275109170305677
$ws.Range("H17").Value = 0.1768558951965065
$ws.Range("I17").Value = 0.08733624454148471
$ws.Range("J17").Value = 0.3842794759825328
$ws.Range("K17").Value = 0.07641921397379912
$ws.Range("M17").Value = 0.01746724890829694
$ws.Range("O17").Value = 0.07205240174672489
$ws.Range("S17").Value = 0.1528384279475982
$ws.Range("F18").Value = 0.03883495145631068
$ws.Range("H18").Value = 0.1310679611650485
$ws.Range("I18").Value = 0.09223300970873786
$ws.Range("J18").Value = 0.4271844660194175
$ws.Range("K18").Value = 0.09223300970873786
$ws.Range("M18").Value = 0.009708737864077669
$ws.Range("O18").Value = 0.07281553398058252
$ws.Range("S18").Value = 0.1359223300970874
$ws.Range("F19").Value = 0.02370030581039755
$ws.Range("H19").Value = 0.1957186544342508
$ws.Range("I19").Value = 0.08486238532110092
$ws.Range("J19").Value = 0.3555045871559633
$ws.Range("K19").Value = 0.1024464831804281
$ws.Range("M19").Value = 0.02140672782874618
$ws.Range("N19").Value = 0.001529051987767584
$ws.Range("O19").Value = 0.077217125382263
$ws.Range("S19").Value = 0.1376146788990826
